$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("With Guidance Sheet")
$ws.Activate()

# Insert two new rows before row 9, pushing the old rows 9-14 down to 11-16
$ws.Rows.Item(9).Resize(2).Insert()

# Update row 8 label: "redStorage LOW" -> "redStorage space1"
$ws.Range("B8").Value = "redStorage space1"

# Fill in the new labels first, then the new Cartesian values, to mirror
# the shared-string insertion order used by the original edit
$ws.Range("B9").Value = "redStorage space2"
$ws.Range("B10").Value = "redStorage space3"
$ws.Range("A9").Value = "Cartesian: 82.73, -228.93, -22.06, -99.67"
$ws.Range("A10").Value = "Cartesian: 85.73, -228.93, 10, -99.67"

# Update the sheet's active-cell selection
$ws.Range("A6").Select()
